$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-31 down to 6-32
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 45169
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 100112017
$ws.Range("G5").Value = "Ramas de apio"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 4000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 4600
$ws.Range("N5").Value = "$/paquete"
$ws.Range("O5").Value = "Región de La Araucanía"
$ws.Range("P5").Value = 4600
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
